# Update market_health_data.xlsx per commit "Update data: 2025-10-30 01:07"

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# 1) Metadata sheet: bump "Last Updated" timestamp by one minute
# -------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "30 Oct 2025, 01:07 AM"

# -------------------------------------------------------------------
# 2) Top Gainers sheet: a new row (SKMEGGPROD) enters the top of the
#    table at row 36, pushing every existing row (36-76) down by one
#    and dropping the former last row (old row 76 / CGPOWER) off the
#    bottom of the 76-row table.
# -------------------------------------------------------------------
$wsGainers = $wb.Worksheets.Item("Top Gainers")

# Insert a blank row at 36 (shifts old rows 36..76 -> 37..77)
$wsGainers.Rows.Item(36).Insert()
# Drop the row that got pushed out past the bottom of the table (old row 76)
$wsGainers.Rows.Item(77).Delete()

# Populate the new row 36 with the incoming data
$wsGainers.Cells.Item(36, 1).Value = "🚀"
$wsGainers.Cells.Item(36, 2).Value = "SKMEGGPROD"
$wsGainers.Cells.Item(36, 3).Value = 4.9959
$wsGainers.Cells.Item(36, 4).Value = 6.6906
$wsGainers.Cells.Item(36, 5).Value = 23.7638

# -------------------------------------------------------------------
# 3) Top Losers sheet: refreshed "Weekly" (column D) figures for a
#    handful of rows whose "Latest"/"Monthly" values stayed the same.
# -------------------------------------------------------------------
$wsLosers = $wb.Worksheets.Item("Top Losers")

$wsLosers.Cells.Item(18, 4).Value = -0.062
$wsLosers.Cells.Item(48, 4).Value = 0.05
$wsLosers.Cells.Item(56, 4).Value = 3.7771
